# Slide 7 ("The Scope Level of a Variable Declaration (continued)") -
# the code sample in the Content Placeholder is updated:
#   - parameter name "initialValue" is shortened to "initializer"
#     (and the surrounding runs are merged into a single run)
#   - the statement's trailing ";" becomes ");" (closing paren added)
#   - the placeholder shape is narrowed slightly to fit the new text

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(7)
$sh = $s.Shapes.Item(2)

# Narrow the content placeholder: cx 8321040 -> 8229600 EMU (cy unchanged).
$sh.Width = 8229600 / 12700

$tr = $sh.TextFrame.TextRange

# Change the trailing ";" to ");" first (higher offset, so it is not
# affected by the length change of the edit below).
$tr.Characters(199, 1).Text = ");"

# Replace ", initialValue," with ", initializer,".
$tr.Characters(137, 15).Text = ", initializer,"
